# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row for 64764041-408a-4257-becb-b605a2cf66c2.md
# "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G3").Value = "2016-09-01 20:52:53"

# zh-cn sheet: row for 64764041-408a-4257-becb-b605a2cf66c2 file
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
$wsZhCn.Range("H3").Value = "2016-09-01 20:52:48"
$wsZhCn.Range("K3").Value = "2016-09-01 20:53:35"

# de-de sheet: row for 64764041-408a-4257-becb-b605a2cf66c2 file
# "Correspond Handback DateTime" (K)
$wsDeDe.Range("K3").Value = "2016-09-01 20:53:43"
